$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 values are swapped for A and G columns; C2/C3 both updated.
$ws.Range("A2").Value = "A 36578-2022"
$ws.Range("A3").Value = "A 36523-2022"

$ws.Range("C2").Value = 46081
$ws.Range("C3").Value = 46081

$ws.Range("G2").Value = 0.3
$ws.Range("G3").Value = 0.2
